# Add a "Save" column (H) to the s_vals sheet:
#  - H1: header label "Save", styled like the other header cells (B1:G1)
#  - H2: numeric value 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (G1) onto the new
# header cell H1, so "Save" picks up the same bold/border/centered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Header text for the new column
$ws.Range("H1").Value = "Save"

# New data cell (plain numeric, no special style - like the other data cells)
$ws.Range("H2").Value = 1
